# Applies the betexplorer data update for 2023/poland_division-2_2023-2024.xlsx
# - Rotates the match data (columns F:V) among rows 148-150
# - Appends a new row 151 with a new match (Pogon Siedlce vs Wisla Pulawy)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture current (pre-edit) match data for rows 148-150, columns F:V ---
$f148 = $ws.Range("F148").Value()
$g148 = $ws.Range("G148").Value()
$h148 = $ws.Range("H148").Value()
$i148 = $ws.Range("I148").Value()
$j148 = $ws.Range("J148").Value()
$k148 = $ws.Range("K148").Value()
$l148 = $ws.Range("L148").Value()
$m148 = $ws.Range("M148").Value()
$n148 = $ws.Range("N148").Value()
$o148 = $ws.Range("O148").Value()
$p148 = $ws.Range("P148").Value()
$q148 = $ws.Range("Q148").Value()
$r148 = $ws.Range("R148").Value()
$s148 = $ws.Range("S148").Value()
$t148 = $ws.Range("T148").Value()
$u148 = $ws.Range("U148").Value()
$v148 = $ws.Range("V148").Value()

$f149 = $ws.Range("F149").Value()
$g149 = $ws.Range("G149").Value()
$h149 = $ws.Range("H149").Value()
$i149 = $ws.Range("I149").Value()
$j149 = $ws.Range("J149").Value()
$k149 = $ws.Range("K149").Value()
$l149 = $ws.Range("L149").Value()
$m149 = $ws.Range("M149").Value()
$n149 = $ws.Range("N149").Value()
$o149 = $ws.Range("O149").Value()
$p149 = $ws.Range("P149").Value()
$q149 = $ws.Range("Q149").Value()
$r149 = $ws.Range("R149").Value()
$s149 = $ws.Range("S149").Value()
$t149 = $ws.Range("T149").Value()
$u149 = $ws.Range("U149").Value()
$v149 = $ws.Range("V149").Value()

$f150 = $ws.Range("F150").Value()
$g150 = $ws.Range("G150").Value()
$h150 = $ws.Range("H150").Value()
$i150 = $ws.Range("I150").Value()
$j150 = $ws.Range("J150").Value()
$k150 = $ws.Range("K150").Value()
$l150 = $ws.Range("L150").Value()
$m150 = $ws.Range("M150").Value()
$n150 = $ws.Range("N150").Value()
$o150 = $ws.Range("O150").Value()
$p150 = $ws.Range("P150").Value()
$q150 = $ws.Range("Q150").Value()
$r150 = $ws.Range("R150").Value()
$s150 = $ws.Range("S150").Value()
$t150 = $ws.Range("T150").Value()
$u150 = $ws.Range("U150").Value()
$v150 = $ws.Range("V150").Value()

# --- Step 2: write the rotated data back ---
# New row 148 <- old row 150 data
$ws.Range("F148").Value = $f150
$ws.Range("G148").Value = $g150
$ws.Range("H148").Value = $h150
$ws.Range("I148").Value = $i150
$ws.Range("J148").Value = $j150
$ws.Range("K148").Value = $k150
$ws.Range("L148").Value = $l150
$ws.Range("M148").Value = $m150
$ws.Range("N148").Value = $n150
$ws.Range("O148").Value = $o150
$ws.Range("P148").Value = $p150
$ws.Range("Q148").Value = $q150
$ws.Range("R148").Value = $r150
$ws.Range("S148").Value = $s150
$ws.Range("T148").Value = $t150
$ws.Range("U148").Value = $u150
$ws.Range("V148").Value = $v150

# New row 149 <- old row 148 data
$ws.Range("F149").Value = $f148
$ws.Range("G149").Value = $g148
$ws.Range("H149").Value = $h148
$ws.Range("I149").Value = $i148
$ws.Range("J149").Value = $j148
$ws.Range("K149").Value = $k148
$ws.Range("L149").Value = $l148
$ws.Range("M149").Value = $m148
$ws.Range("N149").Value = $n148
$ws.Range("O149").Value = $o148
$ws.Range("P149").Value = $p148
$ws.Range("Q149").Value = $q148
$ws.Range("R149").Value = $r148
$ws.Range("S149").Value = $s148
$ws.Range("T149").Value = $t148
$ws.Range("U149").Value = $u148
$ws.Range("V149").Value = $v148

# New row 150 <- old row 149 data
$ws.Range("F150").Value = $f149
$ws.Range("G150").Value = $g149
$ws.Range("H150").Value = $h149
$ws.Range("I150").Value = $i149
$ws.Range("J150").Value = $j149
$ws.Range("K150").Value = $k149
$ws.Range("L150").Value = $l149
$ws.Range("M150").Value = $m149
$ws.Range("N150").Value = $n149
$ws.Range("O150").Value = $o149
$ws.Range("P150").Value = $p149
$ws.Range("Q150").Value = $q149
$ws.Range("R150").Value = $r149
$ws.Range("S150").Value = $s149
$ws.Range("T150").Value = $t149
$ws.Range("U150").Value = $u149
$ws.Range("V150").Value = $v149

# --- Step 3: append new row 151 ---
# Copy formatting (styles) of columns A:E from row 150 down to row 151,
# matching the style pattern used for other data rows (bold index in A, datetime format in E).
$ws.Range("A150:E150").Copy() | Out-Null
$ws.Range("A151:E151").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A151").Value = 150
$ws.Range("B151").Value = "poland"
$ws.Range("C151").Value = "division-2"
$ws.Range("D151").Value = "2023-2024"
$ws.Range("E151").Value = 45242.83333333334
$ws.Range("F151").Value = "Pogon Siedlce"
$ws.Range("G151").Value = 2
$ws.Range("H151").Value = "Wisla Pulawy"
$ws.Range("I151").Value = 2
$ws.Range("J151").Value = 1.92
$ws.Range("K151").Value = "11/11/2023 08:12"
$ws.Range("L151").Value = 1.72
$ws.Range("M151").Value = "12/11/2023 19:50"
$ws.Range("N151").Value = 3.4
$ws.Range("O151").Value = "11/11/2023 08:12"
$ws.Range("P151").Value = 3.89
$ws.Range("Q151").Value = "12/11/2023 19:50"
$ws.Range("R151").Value = 3.34
$ws.Range("S151").Value = "11/11/2023 08:12"
$ws.Range("T151").Value = 4.25
$ws.Range("U151").Value = "12/11/2023 19:50"
$ws.Range("V151").Value = "https://www.betexplorer.com/football/poland/division-2/pogon-siedlce-wisla-pulawy/fTCSRlZr/"
